# Apply changes described in the commit "Second to last runs"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename: frac_driving -> frac_sailing
$ws.Range("K1").Value = "frac_sailing"

# Row 2
$ws.Range("E2").Value = 250
$ws.Range("F2").Value = 689.9447207751965
$ws.Range("G2").Value = 411.1575575524141
$ws.Range("H2").Value = 78.21534042754499
$ws.Range("I2").Value = 0.0786559674099861
$ws.Range("J2").Value = 0.2394453849136681
$ws.Range("K2").Value = 0.6818986476763458
$ws.Range("L2").Value = 8194.50252070112
$ws.Range("M2").Value = "0        50`n1       111`n2        56`n3        73`n4        50`n       ... `n1460     55`n1461     39`n1462     74`n1463     43`n1464     20`nName: total_charged, Length: 147639, dtype: int64"

# Row 3
$ws.Range("E3").Value = 246
$ws.Range("F3").Value = 688.0886383578821
$ws.Range("G3").Value = 1151.069482556251
$ws.Range("H3").Value = 376.5177835322399
$ws.Range("I3").Value = 0.1940679583731451
$ws.Range("J3").Value = 0.3964474108130775
$ws.Range("K3").Value = 0.4094846308137775
$ws.Range("L3").Value = 8290.212922329702
$ws.Range("M3").Value = "0       440`n1       837`n2       233`n3         0`n4       598`n       ... `n1244    888`n1245     51`n1246    371`n1247    325`n1248    205`nName: total_charged, Length: 126831, dtype: int64"

# Row 4
$ws.Range("E4").Value = 247
$ws.Range("F4").Value = 681.9482143690592
$ws.Range("G4").Value = 797.2282901403619
$ws.Range("H4").Value = 228.669615049932
$ws.Range("I4").Value = 0.1501816188972414
$ws.Range("J4").Value = 0.3490308216444599
$ws.Range("K4").Value = 0.5007875594582987
$ws.Range("L4").Value = 8256.743845019051
$ws.Range("M4").Value = "0         0`n1         0`n2         0`n3       353`n4       302`n       ... `n1312    310`n1313     32`n1314      0`n1315     89`n1316    126`nName: total_charged, Length: 135758, dtype: int64"

# Row 5
$ws.Range("E5").Value = 265
$ws.Range("F5").Value = 735.1793937606573
$ws.Range("G5").Value = 863.9383246671267
$ws.Range("H5").Value = 90.14469634769264
$ws.Range("I5").Value = 0.07538108026612521
$ws.Range("J5").Value = 0.3330091836601908
$ws.Range("K5").Value = 0.5916097360736841
$ws.Range("L5").Value = 11330.76436068804
$ws.Range("M5").Value = "0         6`n1         8`n2         0`n3        77`n4        31`n       ... `n1355    203`n1356    112`n1357    147`n1358     20`n1359    212`nName: total_charged, Length: 140894, dtype: int64"

# Row 6
$ws.Range("E6").Value = 267
$ws.Range("F6").Value = 740.8945746071652
$ws.Range("G6").Value = 1195.444624652926
$ws.Range("H6").Value = 420.0686714340791
$ws.Range("I6").Value = 0.1943338307081175
$ws.Range("J6").Value = 0.3984978333165458
$ws.Range("K6").Value = 0.4071683359753367
$ws.Range("L6").Value = 11242.40371892212
$ws.Range("M6").Value = "0       344`n1       661`n2       643`n3       121`n4       373`n       ... `n1329      0`n1330    276`n1331    225`n1332    114`n1333    564`nName: total_charged, Length: 132117, dtype: int64"

# Row 7
$ws.Range("E7").Value = 269
$ws.Range("F7").Value = 744.8374683771093
$ws.Range("G7").Value = 1096.403923527824
$ws.Range("H7").Value = 258.1790175287239
$ws.Range("I7").Value = 0.1441430084568653
$ws.Range("J7").Value = 0.391017166825935
$ws.Range("K7").Value = 0.4648398247171997
$ws.Range("L7").Value = 11222.12660708321
$ws.Range("M7").Value = "0         0`n1        91`n2         0`n3       296`n4         0`n       ... `n1408     41`n1409      0`n1410      0`n1411    175`n1412    300`nName: total_charged, Length: 136793, dtype: int64"

# Row 8
$ws.Range("E8").Value = 181
$ws.Range("F8").Value = 513.5685510490364
$ws.Range("G8").Value = 248.5731362189647
$ws.Range("H8").Value = 58.64620422562893
$ws.Range("I8").Value = 0.0838175065325907
$ws.Range("J8").Value = 0.151822165247239
$ws.Range("K8").Value = 0.7643603282201704
$ws.Range("L8").Value = 5468.09356784078
$ws.Range("M8").Value = "0        30`n1        56`n2        90`n3        29`n4        55`n       ... `n1279    113`n1280    109`n1281     56`n1282     64`n1283     41`nName: total_charged, Length: 126562, dtype: int64"

# Row 9
$ws.Range("E9").Value = 182
$ws.Range("F9").Value = 510.8959575445595
$ws.Range("G9").Value = 825.86044494143
$ws.Range("H9").Value = 282.3527500613303
$ws.Range("I9").Value = 0.2082898776981312
$ws.Range("J9").Value = 0.3460109130145631
$ws.Range("K9").Value = 0.4456992092873058
$ws.Range("L9").Value = 5433.788536721603
$ws.Range("M9").Value = "0       117`n1       258`n2       346`n3       602`n4       147`n       ... `n1142    111`n1143    319`n1144      0`n1145    601`n1146    480`nName: total_charged, Length: 112511, dtype: int64"

# Row 10
$ws.Range("E10").Value = 180
$ws.Range("F10").Value = 510.0106393677598
$ws.Range("G10").Value = 705.9779011012704
$ws.Range("H10").Value = 170.0962818201573
$ws.Range("I10").Value = 0.1507240000854073
$ws.Range("J10").Value = 0.3279090607207886
$ws.Range("K10").Value = 0.5213669391938041
$ws.Range("L10").Value = 5473.387576327632
$ws.Range("M10").Value = "0       140`n1       319`n2        27`n3       453`n4       326`n       ... `n1194    143`n1195     72`n1196    109`n1197     86`n1198    115`nName: total_charged, Length: 119046, dtype: int64"

